# epitweetr topics.xlsx - add "Monkeypox" topic row
# (Typos and updated users/topics files)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 86 (1-based, "#" = 85) is the last existing data row.
# Copy its formatting down into the brand-new row 87 ("#" = 86) before
# filling in the new topic's values, so the new row looks consistent
# with the rest of the table.
$ws.Range("A86:I86").Copy()
$ws.Range("A87:I87").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(87, 1).Value = 86
$ws.Cells.Item(87, 2).Value = "Monkeypox"
$ws.Cells.Item(87, 3).Value = "Monkeypox"
$ws.Cells.Item(87, 4).Value = 0.025
$ws.Cells.Item(87, 5).Value = 0.05
$ws.Cells.Item(87, 6).Value = 'monkeypox OR "viruela del mono" OR "viruela del simio" OR "variole du singe" OR "variola des macacos"'
$ws.Cells.Item(87, 7).Formula = "=LEN(F87)"
$ws.Cells.Item(87, 8).Formula = '=LEN(TRIM(F87))-LEN(SUBSTITUTE(F87," ",""))+1'
$ws.Cells.Item(87, 9).Value = 1

# Extend the conditional formatting that covered rows 2:86 (and 1:86 for
# column H) so that it also covers the newly added row 87.
$fcs = $ws.Cells.FormatConditions
$fcs.Item(1).ModifyAppliesToRange($ws.Range("G2:G87"))
$fcs.Item(2).ModifyAppliesToRange($ws.Range("H2:H87"))
$fcs.Item(5).ModifyAppliesToRange($ws.Range("H1:H87"))
$fcs.Item(7).ModifyAppliesToRange($ws.Range("G2:G87"))

# Move the view / selection down to the row right after the new entry,
# matching where the author ended up after typing the new row.
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 76
$win.ScrollColumn = 1
$ws.Range("A88").Select() | Out-Null
